$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 27 de Abril de 2020 a las 23:22"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1004942
$ws.Range("C4").Value = 17782
$ws.Range("E4").Value = 810824
$ws.Range("G4").Value = 1114
$ws.Range("H4").Value = 56527

# Row 14 - Brasil
$ws.Range("D14").Value = 31142
$ws.Range("E14").Value = 30816

# Row 15 - Canada
$ws.Range("B15").Value = 48242
$ws.Range("C15").Value = 1347
$ws.Range("D15").Value = 18100
$ws.Range("E15").Value = 27440
$ws.Range("G15").Value = 142
$ws.Range("H15").Value = 2702

# Row 159 - Guyana
$ws.Range("D159").Value = 15
$ws.Range("E159").Value = 51

# Row 170 - Siria
$ws.Range("D170").Value = 19
$ws.Range("E170").Value = 21
